$wb = $excel.ActiveWorkbook
$wsTags = $wb.Worksheets.Item("XML Tags")
$wsIds = $wb.Worksheets.Item("Weakness IDs")

# Insert a new row before row 4 (pushes Weakness_ID_1/2/3 rows down)
$wsTags.Rows.Item(4).Insert()

$wsTags.Range("A4").Value = "Line_Number"
$wsTags.Range("B4").Value = "AnalysisInfo/Unified/Context/FunctionDeclarationSourceLocation/line"
$wsTags.Range("C4").Value = "Attribute"

# Match formatting used elsewhere in the table for this row
$wsTags.Range("A3").Copy()
$wsTags.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTags.Range("B4").Borders.LineStyle = 1
$wsTags.Range("B4").Borders.Weight = 2

$wsIds.Range("B5").Copy()
$wsTags.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsScan = $wb.Worksheets.Item("Scan Commands")

$wsIds.Activate()
$excel.ActiveWindow.Zoom = 100
$wsIds.Range("A15").Select()
$excel.ActiveWindow.ScrollRow = 15

$wsScan.Activate()
$excel.ActiveWindow.Zoom = 100

$wsTags.Activate()
$excel.ActiveWindow.Zoom = 100
$wsTags.Range("C12").Select()
